$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '39.646.04'
$ws.Range('E2').Value = '  +2.25%  '

# Row 3
$ws.Range('D3').Value = '2.157.78'
$ws.Range('E3').Value = '  +2.70%  '

# Row 4
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.01'
$ws.Range('E5').Value = '  -0.14%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.621'
$ws.Range('E6').Value = '  +0.91%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '62.75'
$ws.Range('E7').Value = '  +1.56%  '

# Row 8
$ws.Range('E8').Value = '  +0.08%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.390'
$ws.Range('E9').Value = '  +0.55%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0842'
$ws.Range('E10').Value = '  +0.05%  '

# Row 11
$ws.Range('E11').Value = '  +0.20%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.87'
$ws.Range('E12').Value = '  +0.65%  '

# Row 13
$ws.Range('D13').Value = '2.479.92'
$ws.Range('E13').Value = '  +2.86%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.74'
$ws.Range('E14').Value = '  -1.20%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.804'
$ws.Range('E15').Value = '  +0.51%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.48'
$ws.Range('E16').Value = '  +0.03%  '

# Row 17
$ws.Range('D17').Value = '2.170.12'
$ws.Range('E17').Value = '  +3.91%  '

# Row 18
$ws.Range('D18').Value = '39.631.39'
$ws.Range('E18').Value = '  +2.27%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '71.68'
$ws.Range('E19').Value = '  +0.13%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.02'
$ws.Range('E20').Value = '  +0.01%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0847'
$ws.Range('E21').Value = '  +0.60%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '227.51'
$ws.Range('E22').Value = '  +0.28%  '

# Row 23
$ws.Range('E23').Value = '  +0.03%  '

# Row 24
$ws.Range('E24').Value = '  +1.78%  '

# Row 25
$ws.Range('E25').Value = '  -0.29%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '170.69'
$ws.Range('E26').Value = '  +0.40%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.41'
$ws.Range('E27').Value = '  -2.50%  '

# Row 28
$ws.Range('E28').Value = '  +2.20%  '

# Row 29
$ws.Range('E29').Value = '  +1.66%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '19.60'
$ws.Range('E30').Value = '  +1.38%  '

# Row 31
$ws.Range('E31').Value = '  +5.87%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.121'
$ws.Range('E32').Value = '  +0.83%  '

# Row 33
$ws.Range('E33').Value = '  +0.50%  '

# Row 34
$ws.Range('E34').Value = '  -1.84%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.95'
$ws.Range('E35').Value = '  -3.03%  '

# Row 36
$ws.Range('E36').Value = '  +0.38%  '

# Row 37
$ws.Range('E37').Value = '  +6.93%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.39'
$ws.Range('E38').Value = '  +1.43%  '

# Row 39
$ws.Range('E39').Value = '  -0.14%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.90'
$ws.Range('E40').Value = '  +17.50%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '102.46'
$ws.Range('E41').Value = '  +0.68%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0226'
$ws.Range('E42').Value = '  -1.26%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '17.59'
$ws.Range('E43').Value = '  -2.16%  '

# Row 44
$ws.Range('D44').Value = '1.515.05'
$ws.Range('E44').Value = '  -0.68%  '

# Row 45
$ws.Range('E45').Value = '  +0.30%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.88'
$ws.Range('E46').Value = '  +1.65%  '

# Row 47
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.81'
$ws.Range('E47').Value = '  +0.07%  '

# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0920'
$ws.Range('E48').Value = '  +0.97%  '

# Row 49
$ws.Range('E49').Value = '  +0.74%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '49.43'
$ws.Range('E50').Value = '  +7.41%  '

# Row 51
$ws.Range('B51').Value = 'TerraClassic'
$ws.Range('C51').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.000190'
$ws.Range('E51').Value = '  +27.80%  '
